$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D6").Value = -7.714999999999996
$ws.Range("C7").Value = -12.8578
$ws.Range("E7").Value = 16.03730000000001
$ws.Range("B9").Value = 6.590199999999992
$ws.Range("E10").Value = 16.1668
$ws.Range("C12").Value = -10.8268
$ws.Range("E13").Value = 16.66750000000001
$ws.Range("C14").Value = -13.2666
$ws.Range("D15").Value = -8.929099999999995
$ws.Range("E16").Value = 16.38010000000001
$ws.Range("B18").Value = 5.646599999999995
$ws.Range("B20").Value = 9.077500000000001
$ws.Range("E20").Value = 15.70979999999999
$ws.Range("E24").Value = 16.72540000000001
$ws.Range("C26").Value = -12.5959
$ws.Range("B27").Value = 6.119700000000005
$ws.Range("C27").Value = -12.61569999999999
$ws.Range("C29").Value = -11.22760000000001
$ws.Range("D33").Value = -7.588699999999997
$ws.Range("B35").Value = 8.525900000000005
$ws.Range("D35").Value = -8.641999999999994
$ws.Range("C37").Value = -13.2355
$ws.Range("C38").Value = -12.4748
$ws.Range("D38").Value = -8.802299999999994
$ws.Range("E39").Value = 16.1493
$ws.Range("D43").Value = -8.359300000000003
$ws.Range("D44").Value = -7.558599999999999
$ws.Range("D47").Value = -7.8679
$ws.Range("E47").Value = 16.7112
$ws.Range("E48").Value = 17.31740000000001
$ws.Range("C51").Value = -12.1671
$ws.Range("D51").Value = -7.766200000000001
$ws.Range("C52").Value = -11.3901
$ws.Range("E52").Value = 17.4506
$ws.Range("C55").Value = -13.68399999999999
$ws.Range("E56").Value = 16.64470000000001
$ws.Range("D57").Value = -8.281099999999999
$ws.Range("D63").Value = -7.936299999999997
$ws.Range("B69").Value = 6.377699999999993
$ws.Range("C69").Value = -11.4127
$ws.Range("C70").Value = -12.4845
$ws.Range("D70").Value = -8.132
$ws.Range("B76").Value = 4.464699999999999
$ws.Range("B78").Value = 9.648900000000003
$ws.Range("C81").Value = -12.7585
$ws.Range("B82").Value = 5.702399999999999
$ws.Range("B83").Value = 5.960599999999996
$ws.Range("C83").Value = -14.00840000000001
$ws.Range("E84").Value = 17.089
$ws.Range("D88").Value = -7.277899999999996
$ws.Range("B93").Value = 6.468399999999997
$ws.Range("D99").Value = -7.712999999999997
$ws.Range("E100").Value = 16.44000000000001
$ws.Range("E101").Value = 16.80250000000001
$ws.Range("C102").Value = -13.3546